$d = $word.ActiveDocument
$xml = @'
<w:tbl><w:tblGrid><w:gridCol w:w="3515" w:type="dxa"/><w:gridCol w:w="8391" w:type="dxa"/></w:tblGrid><w:tblPr><w:tblW w:w="11906" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblpPr w:vertAnchor="page" w:horzAnchor="page" w:tblpXSpec="center" w:tblpYSpec="bottom"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="5386" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3515" w:type="dxa"/></w:tcPr><w:tbl><w:tblGrid><w:gridCol/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/><w:tblCellMar><w:top w:w="283" w:type="dxa"/><w:left w:w="283" w:type="dxa"/><w:right w:w="283" w:type="dxa"/><w:bottom w:w="283" w:type="dxa"/></w:tblCellMar></w:tblPr><w:tr><w:trPr/><w:tc><w:tcPr/><w:tbl><w:tblGrid><w:gridCol/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="397" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr/><w:r><w:rPr><w:rStyle w:val="SwissBill Title"/></w:rPr><w:t xml:space="preserve">Receipt</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="3175" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Heading"/></w:rPr><w:t xml:space="preserve">Account / Payable to</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve">CH93 0076 2011 6238 5295 7</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve">Thomas LeClaire</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve">Rue examplaire 22a</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve">1000 Lausanne</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Heading"/></w:rPr><w:t xml:space="preserve">Reference</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve">RF18 5390 0754 7034</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Value"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Heading"/></w:rPr><w:t xml:space="preserve">Payable by (name/address)</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:pict><v:shape type="#_x0000_t75" stroked="f" style="width:147pt; height:57pt; margin-left:0pt; margin-top:0pt; position:relative; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:margin; mso-position-vertical-relative:line; z-index:2147483647;"><v:imagedata r:id="rId7" o:title=""/></v:shape></w:pict></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="794" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:tbl><w:tblGrid><w:gridCol w:w="692" w:type="dxa"/><w:gridCol w:w="2256" w:type="dxa"/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="794" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="692" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Heading"/></w:rPr><w:t xml:space="preserve">Currency</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Amount"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Amount"/></w:rPr><w:t xml:space="preserve">CHF</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2256" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Heading"/></w:rPr><w:t xml:space="preserve">Amount</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Receipt Amount"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Receipt Amount"/></w:rPr><w:t xml:space="preserve">25.90</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="1020" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr><w:pStyle w:val="SwissBill Acceptance point"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Acceptance point"/></w:rPr><w:t xml:space="preserve">Acceptance point</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc></w:tr></w:tbl><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="8391" w:type="dxa"/></w:tcPr><w:tbl><w:tblGrid><w:gridCol/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/><w:tblCellMar><w:top w:w="283" w:type="dxa"/><w:left w:w="283" w:type="dxa"/><w:right w:w="283" w:type="dxa"/><w:bottom w:w="283" w:type="dxa"/></w:tblCellMar></w:tblPr><w:tr><w:trPr/><w:tc><w:tcPr/><w:tbl><w:tblGrid><w:gridCol w:w="2891" w:type="dxa"/><w:gridCol w:w="4932" w:type="dxa"/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="4819" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2891" w:type="dxa"/></w:tcPr><w:tbl><w:tblGrid><w:gridCol/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="397" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr/><w:r><w:rPr><w:rStyle w:val="SwissBill Title"/></w:rPr><w:t xml:space="preserve">Payment part</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="283" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:rPr><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="2891" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:p><w:pPr/><w:r><w:pict><v:shape type="#_x0000_t75" stroked="f" style="width:130pt; height:130pt; margin-left:0pt; margin-top:0pt; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line;"><w10:wrap type="inline"/><v:imagedata r:id="rId8" o:title=""/></v:shape></w:pict></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="1247" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr/><w:tbl><w:tblGrid><w:gridCol w:w="816" w:type="dxa"/><w:gridCol w:w="2075" w:type="dxa"/></w:tblGrid><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:bidiVisual w:val="0"/></w:tblPr><w:tr><w:trPr><w:trHeight w:val="1247" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="816" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Heading"/></w:rPr><w:t xml:space="preserve">Currency</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Amount"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Amount"/></w:rPr><w:t xml:space="preserve">CHF</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2075" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Heading"/></w:rPr><w:t xml:space="preserve">Amount</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Amount"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Amount"/></w:rPr><w:t xml:space="preserve">25.90</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc></w:tr></w:tbl><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="4932" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Heading"/></w:rPr><w:t xml:space="preserve">Account / Payable to</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve">CH93 0076 2011 6238 5295 7</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve">Thomas LeClaire</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve">Rue examplaire 22a</w:t></w:r><w:br/><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve">1000 Lausanne</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Heading"/></w:rPr><w:t xml:space="preserve">Reference</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve">RF18 5390 0754 7034</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Value"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Value"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Heading"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Heading"/></w:rPr><w:t xml:space="preserve">Payable by (name/address)</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:pict><v:shape type="#_x0000_t75" stroked="f" style="width:184pt; height:71pt; margin-left:0pt; margin-top:0pt; position:relative; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:margin; mso-position-vertical-relative:line; z-index:2147483647;"><v:imagedata r:id="rId9" o:title=""/></v:shape></w:pict></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="567" w:hRule="atLeast"/></w:trPr><w:tc><w:tcPr><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="SwissBill Payment part Further information"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SwissBill Payment part Further information"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc></w:tr></w:tbl><w:p/></w:tc></w:tr></w:tbl>
'@
$d.Content.InsertXML($xml)
